$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(1.063631372210698, 0.6069138574468314, 2.188368220279388, 3.8005269855140895)"
$ws.Range("C2").Value = "NIG(0.9192664996201568, 0.6695023852587427, 5.654705693043235, 5.007920253508492)"
$ws.Range("D2").Value = "NIG(1.8222570068211121, 1.222534893275363, 0.9850379134693319, 4.297237844550496)"
$ws.Range("E2").Value = "JSU(-1.240442123002143, 1.1575792410004593, 2.3695611535422665, 4.825346040271131)"
